$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": a new daily price column ("19-dec") is inserted
#     right before the old "01-oct." block (column EV), shifting every
#     column from EV..FZ one place to the right (EV->EW ... FZ->GA). ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns("EV").Insert()
$ws1.Range("EV1").Value = "19-dec"
for ($r = 2; $r -le 25; $r++) {
    $ws1.Range("EV" + $r).Value = "-"
}

# --- Sheet "Gaz": one more day of data appended as row 182. ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A182").NumberFormat = "@"
$ws2.Range("A182").Value = "2025-12-17"
$ws2.Range("A182").Style = "Normal"
$ws2.Range("B182").Value = 25.75

# --- Sheet "CO2": one more day of data appended as row 182. ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A182").NumberFormat = "@"
$ws3.Range("A182").Value = "2025-12-17"
$ws3.Range("A182").Style = "Normal"
$ws3.Range("B182").Value = 84.8
